# Lookup tables update: remove the obsolete "EQUITAS" symbol row.
#
# The sheet lists trading symbols twice (column A = "Symbol", column B =
# "Table"); the "EQUITAS" entry (row 187, right above "EQUITASBNK") is a
# duplicate/obsolete symbol that the author removed. Deleting the entire
# worksheet row removes both cells, shifts every following row up by one,
# and drops the now-unused "EQUITAS" shared string, which is exactly the
# content change described by the diff (uniqueCount 199 -> 198, dimension
# A1:B191 -> A1:B190).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(187).Delete()

# Reflect the author's resulting selection in the saved view state.
$ws.Range("H182").Select()
